$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.790.24'
$ws.Range('E2').Value = '  +4.22%  '
$ws.Range('D3').Value = '1.928.16'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.68'
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.34'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.76'
$ws.Range('E9').Value = '  +9.15%  '
$ws.Range('E10').Value = '  +3.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0768'
$ws.Range('E11').Value = '  +3.93%  '
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.76'
$ws.Range('E13').Value = '  +8.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.823'
$ws.Range('E14').Value = '  +7.69%  '
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('E16').Value = '  +4.68%  '
$ws.Range('D17').Value = '1.931.00'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '36.695.81'
$ws.Range('E18').Value = '  +3.80%  '
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('D20').Value = '0.0₃0867'
$ws.Range('E20').Value = '  +5.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '252.31'
$ws.Range('E21').Value = '  +3.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.45'
$ws.Range('E22').Value = '  +4.70%  '
$ws.Range('E23').Value = '  +5.55%  '
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.23'
$ws.Range('E26').Value = '  +1.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.36'
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.91'
$ws.Range('E28').Value = '  +4.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.83'
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.58'
$ws.Range('E31').Value = '  +6.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0619'
$ws.Range('E32').Value = '  +4.48%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.97'
$ws.Range('E33').Value = '  -3.68%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.39'
$ws.Range('E34').Value = '  +5.77%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0870'
$ws.Range('E36').Value = '  +20.23%  '
$ws.Range('E37').Value = '  -11.47%  '
$ws.Range('E38').Value = '  +7.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.87'
$ws.Range('E39').Value = '  +50.28%  '
$ws.Range('E40').Value = '  +6.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.60'
$ws.Range('E41').Value = '  +11.86%  '
$ws.Range('E42').Value = '  +5.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.39'
$ws.Range('E43').Value = '  -1.78%  '
$ws.Range('E44').Value = '  +4.15%  '
$ws.Range('D45').Value = '1.345.52'
$ws.Range('E45').Value = '  +3.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.60'
$ws.Range('E46').Value = '  +9.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.81'
$ws.Range('E49').Value = '  +2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.46'
$ws.Range('E50').Value = '  +3.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.58'
$ws.Range('E51').Value = '  +3.72%  '
